$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: date 20.10.18 -> 30.10.18, end time 19:00 -> 18:15, remark loses "-Timer" line
$ws.Range("A18").Value = "30.10.18"
$ws.Range("C18").Value = 0.76041666666666663
$ws.Range("E18").Value = "-Cascasde-Classifier Options and refactoring"

# Row 19: previously-empty placeholder row becomes a real diary entry
$ws.Range("A19").Value = "31.10.18"
$ws.Range("B19").Value = 0.45833333333333331
$ws.Range("C19").Value = 0.75
$ws.Range("E19").Value = "-Cascasde-Classifier Options and refactoring`n-Timer"
$ws.Range("E19").WrapText = $true

# Row heights: 18 shrinks back to single-line default, 19 grows to fit the wrapped remark
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(19).RowHeight = 30

# Selection moves down one row
[void]$ws.Range("C20").Select()
